# "Added last minute updates"
#
# The first paragraph of the document carries a merge-field style marker
# "**ID__AFFARS_pgi_5325_topic_3__ID**" followed by a run containing a
# single trailing space. This change:
#   1. Renames the marker to "**ID__AFFARS_AFMC_PGI_5325_7901_3_90__ID**"
#   2. Removes the now-unwanted trailing-space run
#   3. Gives the paragraph the same paragraph border / indentation
#      treatment already used by the paragraphs further down in the
#      document (w:pBdr space=5 on all sides, w:ind left=225 twips).

$d = $word.ActiveDocument

$oldMarker = "**ID__AFFARS_pgi_5325_topic_3__ID**"
$newMarker = "**ID__AFFARS_AFMC_PGI_5325_7901_3_90__ID**"

# Locate the marker text robustly (rather than assuming a paragraph index).
$matchRange = $d.Content.Duplicate
$found = $matchRange.Find.Execute($oldMarker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # The paragraph that holds the marker.
    $para = $matchRange.Paragraphs(1)

    # The lone space run sits immediately after the marker text; drop it.
    $spaceRange = $d.Range($matchRange.End, $matchRange.End + 1)
    if ($spaceRange.Text -eq " ") {
        $spaceRange.Delete()
    }

    # Swap in the new identifier text.
    $matchRange.Text = $newMarker

    # Match the indentation / paragraph-border formatting used elsewhere
    # in this document (225 twips = 11.25 pt; border space = 5 pt).
    $pf = $para.Range.ParagraphFormat
    $pf.LeftIndent = 11.25
    $pf.Borders.DistanceFromTop = 5
    $pf.Borders.DistanceFromBottom = 5
    $pf.Borders.DistanceFromLeft = 5
    $pf.Borders.DistanceFromRight = 5
}
